# Insert a new weekly record at row 15 (shifting the existing rows 15:96 down to 16:97),
# then populate the new row with the latest Coco price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 15; everything below (15:96) shifts down to (16:97).
$ws.Rows("15:15").Insert()

# Fill in the new row 15 with the new weekly observation.
$ws.Cells.Item(15, 1).Value = 6
$ws.Cells.Item(15, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(15, 3).Value = "Metropolitana"
$ws.Cells.Item(15, 4).Value = 45243
$ws.Cells.Item(15, 5).Value = 13
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100108
$ws.Cells.Item(15, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(15, 9).Value = 100108007
$ws.Cells.Item(15, 10).Value = "Coco"
$ws.Cells.Item(15, 11).Value = "Sin especificar"
$ws.Cells.Item(15, 12).Value = "Primera"
$ws.Cells.Item(15, 13).Value = 86
$ws.Cells.Item(15, 14).Value = 30000
$ws.Cells.Item(15, 15).Value = 30000
$ws.Cells.Item(15, 16).Value = 30000
$ws.Cells.Item(15, 17).Value = "`$/malla 20 unidades"
$ws.Cells.Item(15, 18).Value = "Perú"
$ws.Cells.Item(15, 19).Value = 1500
$ws.Cells.Item(15, 20).Value = 20

# Make sure the new date cell keeps the same date-time number format used by column D elsewhere.
$ws.Cells.Item(15, 4).NumberFormat = $ws.Cells.Item(16, 4).NumberFormat
